# Insert a new weekly record for "Femacal de La Calera - Cebollín" as row 403,
# shifting the existing rows 403:429 down to 404:430.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 403 (pushes old row 403 -> 404, ..., old 429 -> 430)
$ws.Rows.Item(403).Insert()

# Populate the newly inserted row 403 with the new observation
$ws.Range("A403").Value = 3
$ws.Range("B403").Value = "Femacal de La Calera"
$ws.Range("C403").Value = "Coquimbo"
$ws.Range("D403").Value = 44610
$ws.Range("E403").Value = 5
$ws.Range("F403").Value = 100112037
$ws.Range("G403").Value = "Cebollín"
$ws.Range("H403").Value = "Sin especificar"
$ws.Range("I403").Value = "Primera"
$ws.Range("J403").Value = 170
$ws.Range("K403").Value = 3500
$ws.Range("L403").Value = 3800
$ws.Range("M403").Value = 3641
$ws.Range("N403").Value = "$/paquete 36 unidades"
$ws.Range("O403").Value = "Provincia de Quillota"
$ws.Range("P403").Value = 101
$ws.Range("Q403").Value = 36
$ws.Range("R403").Value = "Hortaliza"
